# This workbook tracks per-game "save data" rows (one row per game date).
# Column G holds "K" (renamed from the old "Strike#" label) and has been
# regenerated from an external calc (std/mean + s_vals), so we re-write the
# resulting values cell by cell to match the regenerated data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 1
    4  = 0
    5  = 1
    6  = 1
    7  = 3
    8  = 1
    9  = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 1
    20 = 0
    21 = 2
    22 = 1
    23 = 0
    24 = 2
    25 = 3
    26 = 1
    27 = 2
    29 = 0
    30 = 0
    31 = 2
    32 = 1
    33 = 0
    34 = 0
    35 = 3
    36 = 0
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 2
    42 = 0
    43 = 0
    44 = 0
    46 = 1
    47 = 1
    48 = 1
    49 = 1
    50 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
